# Update the "bj_sch_sth_impact_202401_3_kato_katz" form to v3.3:
#  - Rename the "Résultats" calculate-group ("bj_k2") to "Résultats KK" ("bj_k2_2")
#  - Replace the participant-code duplicate check: swap the old 3-digit/2-digit
#    regex constraint for a "not already used" uniqueness constraint (and message)
#  - Bump the form title / form_id on the settings sheet from V3.2 to V3.3

$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")

# --- survey sheet: row 9 is the "begin repeat" group holding the results calc ---
$wsSurvey.Range("B9").Value = "bj_k2_2"

# --- settings sheet: bump form_title / form_id to V3.3 ---
$wsSettings.Range("A2").Value = "(2024 Janvier) - 3. SCH/STH - Kato Katz V3.3"
$wsSettings.Range("B2").Value = "bj_sch_sth_impact_202401_3_kato_katz_v3_3"

# --- survey sheet: row 14 is the k_espen_code_id field (constraint / message) ---
$wsSurvey.Range("G14").Value = "Cet identifiant est déjà utilisé"
$wsSurvey.Range("F14").Value = 'not(selected(${C3}, ${k_espen_code_id}))'

# --- survey sheet: finish renaming the results group label ---
$wsSurvey.Range("C9").Value = "Résultats KK"

# --- restore the cursor/selection positions left by the edit ---
$wsSettings.Range("F30").Select()
$wsSurvey.Range("C9").Select()
